$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force the Price column (D) to text formatting while writing the
# updated values, so numeric-looking strings (e.g. '2.62') are stored as text
# instead of being auto-converted to numbers -- then restore the original
# (default/no explicit format) style so no stray formatting change is left behind.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '43.374.04'
$ws.Range("E2").Value = '  +2.61%  '
$ws.Range("D3").Value = '2.305.28'
$ws.Range("E3").Value = '  +1.49%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '310.68'
$ws.Range("E5").Value = '  +0.60%  '
$ws.Range("D6").Value = '103.64'
$ws.Range("E6").Value = '  +6.72%  '
$ws.Range("D7").Value = '0.533'
$ws.Range("E7").Value = '  +1.04%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +7.73%  '
$ws.Range("D10").Value = '36.63'
$ws.Range("E10").Value = '  +4.19%  '
$ws.Range("D11").Value = '52.57'
$ws.Range("E11").Value = '  +0.88%  '
$ws.Range("E12").Value = '  +0.36%  '
$ws.Range("E13").Value = '  -1.12%  '
$ws.Range("D14").Value = '7.00'
$ws.Range("E14").Value = '  +3.13%  '
$ws.Range("D15").Value = '2.663.56'
$ws.Range("E15").Value = '  +1.55%  '
$ws.Range("D16").Value = '15.08'
$ws.Range("E16").Value = '  +2.93%  '
$ws.Range("D17").Value = '2.302.87'
$ws.Range("E17").Value = '  +1.29%  '
$ws.Range("D18").Value = '0.809'
$ws.Range("E18").Value = '  +2.32%  '
$ws.Range("D19").Value = '43.263.65'
$ws.Range("E19").Value = '  +2.74%  '
$ws.Range("D20").Value = '12.15'
$ws.Range("E20").Value = '  -0.82%  '
$ws.Range("D21").Value = '0.0₃0926'
$ws.Range("E21").Value = '  +2.18%  '
$ws.Range("D22").Value = '6.16'
$ws.Range("E22").Value = '  +3.08%  '
$ws.Range("D23").Value = '68.08'
$ws.Range("E23").Value = '  +0.55%  '
$ws.Range("D24").Value = '242.57'
$ws.Range("E24").Value = '  +2.62%  '
$ws.Range("E25").Value = '  +2.68%  '
$ws.Range("D26").Value = '2.62'
$ws.Range("E26").Value = '  +1.11%  '
$ws.Range("E27").Value = '  +0.16%  '
$ws.Range("D28").Value = '24.87'
$ws.Range("E28").Value = '  +5.60%  '
$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").Value = '36.80'
$ws.Range("E29").Value = '  +0.38%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '2.22'
$ws.Range("E30").Value = '  +4.13%  '
$ws.Range("D31").Value = '9.64'
$ws.Range("E31").Value = '  +0.90%  '
$ws.Range("D32").Value = '167.59'
$ws.Range("E32").Value = '  +2.09%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("E35").Value = '  +3.85%  '
$ws.Range("E36").Value = '  +6.30%  '
$ws.Range("E37").Value = '  +0.76%  '
$ws.Range("E38").Value = '  -2.17%  '
$ws.Range("E39").Value = '  +2.94%  '
$ws.Range("E40").Value = '  +1.72%  '
$ws.Range("E41").Value = '  +7.29%  '
$ws.Range("D43").Value = '2.65'
$ws.Range("E43").Value = '  +16.41%  '
$ws.Range("D44").Value = '0.0293'
$ws.Range("E44").Value = '  +4.09%  '
$ws.Range("D45").Value = '1.981.26'
$ws.Range("E45").Value = '  +1.42%  '
$ws.Range("D46").Value = '19.05'
$ws.Range("E46").Value = '  +1.21%  '
$ws.Range("E47").Value = '  +3.14%  '
$ws.Range("E48").Value = '  +1.54%  '
$ws.Range("D49").Value = '55.90'
$ws.Range("E49").Value = '  +4.30%  '
$ws.Range("D50").Value = '2.94'
$ws.Range("E50").Value = '  -0.05%  '
$ws.Range("E51").Value = '  +8.23%  '

$ws.Range("D2:D51").Style = "Normal"
